$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Third Choice" entries for Jeremy (row 5) and Jim (row 9)
$ws.Range("D5").ClearContents()
$ws.Range("D9").ClearContents()

# Move the active selection to D9
$ws.Range("D9").Select()
